# filter() method for featureGroupsSet and negate support for features method
#
# Updates the "fGroups" implementation-status matrix on sheet1:
#  - adds a new "set" column (G) marking several methods as implemented ("X")
#  - promotes several "maybe" (X?) entries in column B to definite ("X")
#  - moves a couple of stray marks from one column to the correct one
#  - adds a remark ("maybe wait for autoID branch") in a new column H for
#    two rows that still need further work

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "X" marks added in column G (the new "set" implementation column)
$newG = @(7, 8, 12, 13, 14, 24, 30, 31, 32, 34, 35, 36, 37, 38, 44, 50, 51, 53)
foreach ($r in $newG) {
    $ws.Range("G$r").Value = "X"
}

# Column B entries promoted from "X?" (maybe) to "X" (done)
$promoteB = @(8, 30, 31, 34, 44, 50, 51)
foreach ($r in $promoteB) {
    $ws.Range("B$r").Value = "X"
}

# featureTable row: mark moved from B13 to C13
$ws.Range("B13").ClearContents()
$ws.Range("C13").Value = "X"

# groupInfo row: mark moved from C24 to G24
$ws.Range("C24").ClearContents()
$ws.Range("G24").Value = "X"

# Remarks added in new column H
$ws.Range("H23").Value = "maybe wait for autoID branch"
$ws.Range("H48").Value = "maybe wait for autoID branch"

# Update the active selection to match the edited area
[void]$ws.Range("G15").Select()
